$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.38"
$ws.Range("E2").Value = "'-0.30%"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'31.65"
$ws.Range("E3").Value = "'0.59%"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'5.150"
$ws.Range("E4").Value = "'1.18%"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.08119"
$ws.Range("E5").Value = "'10.46%"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'2.516"
$ws.Range("E6").Value = "'13.96%"
$ws.Range("G6").Value = "'22"
$ws.Range("E7").Value = "'-1.66%"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'3.888"
$ws.Range("E8").Value = "'2.75%"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'0.9248"
$ws.Range("E9").Value = "'0.55%"
$ws.Range("G9").Value = "'22"
$ws.Range("D10").Value = "'0.1754"
$ws.Range("E10").Value = "'2.96%"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.07438"
$ws.Range("E11").Value = "'-1.93%"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'0.08912"
$ws.Range("E12").Value = "'9.41%"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'0.03028"
$ws.Range("E13").Value = "'0.10%"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'0.001515"
$ws.Range("E15").Value = "'1.42%"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.005849"
$ws.Range("E16").Value = "'-4.09%"
$ws.Range("G16").Value = "'22"
$ws.Range("D17").Value = "'3.574"
$ws.Range("E17").Value = "'3.19%"
$ws.Range("G17").Value = "'22"
$ws.Range("D18").Value = "'2.285"
$ws.Range("E18").Value = "'2.70%"
$ws.Range("G18").Value = "'22"
$ws.Range("E19").Value = "'-0.47%"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'0.1340"
$ws.Range("E20").Value = "'0.24%"
$ws.Range("G20").Value = "'22"
$ws.Range("D21").Value = "'4.158"
$ws.Range("E21").Value = "'-10.77%"
$ws.Range("G21").Value = "'22"
$ws.Range("D22").Value = "'0.1682"
$ws.Range("E22").Value = "'7.38%"
$ws.Range("G22").Value = "'22"
$ws.Range("D23").Value = "'0.04620"
$ws.Range("E23").Value = "'-0.52%"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'1.05%"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.004530"
$ws.Range("E25").Value = "'1.11%"
$ws.Range("G25").Value = "'22"
$ws.Range("E26").Value = "'-7.65%"
$ws.Range("G26").Value = "'22"
$ws.Range("D27").Value = "'0.0003412"
$ws.Range("E27").Value = "'27.91%"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("G38").Value = "'22"
$ws.Range("D39").Value = "'0.01764"
$ws.Range("E39").Value = "'1.72%"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.04600"
$ws.Range("E40").Value = "'1.54%"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.006947"
$ws.Range("E41").Value = "'-4.01%"
$ws.Range("G41").Value = "'22"
$ws.Range("E42").Value = "'2.02%"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.002191"
$ws.Range("E43").Value = "'-1.73%"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.01067"
$ws.Range("E44").Value = "'-0.83%"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.00006223"
$ws.Range("E45").Value = "'-1.19%"
$ws.Range("G45").Value = "'22"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("G46").Value = "'22"
$ws.Range("D47").Value = "'0.008404"
$ws.Range("E47").Value = "'-15.89%"
$ws.Range("G47").Value = "'22"
$ws.Range("D48").Value = "'0.7480"
$ws.Range("E48").Value = "'-7.48%"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("G49").Value = "'22"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("G50").Value = "'22"
$ws.Range("G51").Value = "'22"
